$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.41996
$ws.Range("H2").Value = 88.25988000000001
$ws.Range("I2").Value = 0.6207199949605289
$ws.Range("J2").Value = 0.66829493802317
$ws.Range("M2").Value = 26.07194833333334
$ws.Range("N2").Value = 78.215845
$ws.Range("O2").Value = 0.7595928012803946
$ws.Range("P2").Value = 0.7818221335728009
$ws.Range("Q2").Value = 767.0356770887335
$ws.Range("R2").Value = 6903.3210937986
$ws.Range("S2").Value = 0.4714944397828206
$ws.Range("T2").Value = 0.5224877743011775

$ws.Range("G3").Value = 29.41996
$ws.Range("H3").Value = 88.25988000000001
$ws.Range("I3").Value = 0.6207199949605289
$ws.Range("J3").Value = 0.66829493802317
$ws.Range("O3").Value = 0.155109029208254
$ws.Range("P3").Value = 0.1596482641062294
$ws.Range("Q3").Value = 156.62886620408
$ws.Range("R3").Value = 1409.65979583672
$ws.Range("S3").Value = 0.09627927582847998
$ws.Range("T3").Value = 0.1066921267663792

$ws.Range("G4").Value = 29.41996
$ws.Range("H4").Value = 88.25988000000001
$ws.Range("I4").Value = 0.6207199949605289
$ws.Range("J4").Value = 0.66829493802317
$ws.Range("M4").Value = 2.927739
$ws.Range("N4").Value = 5.855478
$ws.Range("O4").Value = 0.08529816951135136
$ws.Range("P4").Value = 0.05852960232096958
$ws.Range("Q4").Value = 86.13396427044
$ws.Range("R4").Value = 516.80378562264
$ws.Range("S4").Value = 0.05294627934922835
$ws.Range("T4").Value = 0.03911503695561316

$ws.Range("I5").Value = 0.1515698101047853
$ws.Range("J5").Value = 0.1631868437822795
$ws.Range("M5").Value = 26.07194833333334
$ws.Range("N5").Value = 78.215845
$ws.Range("O5").Value = 0.7595928012803946
$ws.Range("P5").Value = 0.7818221335728009
$ws.Range("Q5").Value = 187.2977395022172
$ws.Range("R5").Value = 1685.679655519955
$ws.Range("S5").Value = 0.1151313366470313
$ws.Range("T5").Value = 0.1275830863768731

$ws.Range("I6").Value = 0.1515698101047853
$ws.Range("J6").Value = 0.1631868437822795
$ws.Range("O6").Value = 0.155109029208254
$ws.Range("P6").Value = 0.1596482641062294
$ws.Range("S6").Value = 0.02350984610263265
$ws.Range("T6").Value = 0.02605249633481535

$ws.Range("I7").Value = 0.1515698101047853
$ws.Range("J7").Value = 0.1631868437822795
$ws.Range("M7").Value = 2.927739
$ws.Range("N7").Value = 5.855478
$ws.Range("O7").Value = 0.08529816951135136
$ws.Range("P7").Value = 0.05852960232096958
$ws.Range("Q7").Value = 21.032524671407
$ws.Range("R7").Value = 126.195148028442
$ws.Range("S7").Value = 0.01292862735512131
$ws.Range("T7").Value = 0.009551261070591007

$ws.Range("G8").Value = 0.3873096666666667
$ws.Range("H8").Value = 1.161929
$ws.Range("I8").Value = 0.008171692087327698
$ws.Range("J8").Value = 0.008798009571759262
$ws.Range("M8").Value = 26.07194833333334
$ws.Range("N8").Value = 78.215845
$ws.Range("O8").Value = 0.7595928012803946
$ws.Range("P8").Value = 0.7818221335728009
$ws.Range("Q8").Value = 10.09791761833389
$ws.Range("R8").Value = 90.881258565005
$ws.Range("S8").Value = 0.006207158483814082
$ws.Range("T8").Value = 0.006878478614586751

$ws.Range("G9").Value = 0.3873096666666667
$ws.Range("H9").Value = 1.161929
$ws.Range("I9").Value = 0.008171692087327698
$ws.Range("J9").Value = 0.008798009571759262
$ws.Range("O9").Value = 0.155109029208254
$ws.Range("P9").Value = 0.1596482641062294
$ws.Range("Q9").Value = 2.061997159747333
$ws.Range("R9").Value = 18.557974437726
$ws.Range("S9").Value = 0.00126750322665417
$ws.Range("T9").Value = 0.001404586955721357

$ws.Range("G10").Value = 0.3873096666666667
$ws.Range("H10").Value = 1.161929
$ws.Range("I10").Value = 0.008171692087327698
$ws.Range("J10").Value = 0.008798009571759262
$ws.Range("M10").Value = 2.927739
$ws.Range("N10").Value = 5.855478
$ws.Range("O10").Value = 0.08529816951135136
$ws.Range("P10").Value = 0.05852960232096958
$ws.Range("Q10").Value = 1.133941616177
$ws.Range("R10").Value = 6.803649697061999
$ws.Range("S10").Value = 0.0006970303768594466
$ws.Range("T10").Value = 0.0005149440014511536

$ws.Range("G11").Value = 10.122265
$ws.Range("H11").Value = 20.24453
$ws.Range("I11").Value = 0.2135656295858028
$ws.Range("J11").Value = 0.153289545846405
$ws.Range("M11").Value = 26.07194833333334
$ws.Range("N11").Value = 78.215845
$ws.Range("O11").Value = 0.7595928012803946
$ws.Range("P11").Value = 0.7818221335728009
$ws.Range("Q11").Value = 263.9071700963083
$ws.Range("R11").Value = 1583.44302057785
$ws.Range("S11").Value = 0.1622229148342911
$ws.Range("T11").Value = 0.119845159788042

$ws.Range("G12").Value = 10.122265
$ws.Range("H12").Value = 20.24453
$ws.Range("I12").Value = 0.2135656295858028
$ws.Range("J12").Value = 0.153289545846405
$ws.Range("O12").Value = 0.155109029208254
$ws.Range("P12").Value = 0.1596482641062294
$ws.Range("Q12").Value = 53.88990638896999
$ws.Range("R12").Value = 323.3394383338199
$ws.Range("S12").Value = 0.03312595747730345
$ws.Range("T12").Value = 0.02447240990001083

$ws.Range("G13").Value = 10.122265
$ws.Range("H13").Value = 20.24453
$ws.Range("I13").Value = 0.2135656295858028
$ws.Range("J13").Value = 0.153289545846405
$ws.Range("M13").Value = 2.927739
$ws.Range("N13").Value = 5.855478
$ws.Range("O13").Value = 0.08529816951135136
$ws.Range("P13").Value = 0.05852960232096958
$ws.Range("Q13").Value = 29.63535000883499
$ws.Range("R13").Value = 118.54140003534
$ws.Range("S13").Value = 0.01821675727420828
$ws.Range("T13").Value = 0.008971976158352119

$ws.Range("G14").Value = 0.2830933333333334
$ws.Range("H14").Value = 0.84928
$ws.Range("I14").Value = 0.005972873261555284
$ws.Range("J14").Value = 0.006430662776386256
$ws.Range("M14").Value = 26.07194833333334
$ws.Range("N14").Value = 78.215845
$ws.Range("O14").Value = 0.7595928012803946
$ws.Range("P14").Value = 0.7818221335728009
$ws.Range("Q14").Value = 7.380794760177779
$ws.Range("R14").Value = 66.42715284160001
$ws.Range("S14").Value = 0.004536951532437545
$ws.Range("T14").Value = 0.005027634492121494

$ws.Range("G15").Value = 0.2830933333333334
$ws.Range("H15").Value = 0.84928
$ws.Range("I15").Value = 0.005972873261555284
$ws.Range("J15").Value = 0.006430662776386256
$ws.Range("O15").Value = 0.155109029208254
$ws.Range("P15").Value = 0.1596482641062294
$ws.Range("Q15").Value = 1.507160031146667
$ws.Range("R15").Value = 13.56444028032
$ws.Range("S15").Value = 0.000926446573183778
$ws.Range("T15").Value = 0.001026644149302611

$ws.Range("G16").Value = 0.2830933333333334
$ws.Range("H16").Value = 0.84928
$ws.Range("I16").Value = 0.005972873261555284
$ws.Range("J16").Value = 0.006430662776386256
$ws.Range("M16").Value = 2.927739
$ws.Range("N16").Value = 5.855478
$ws.Range("O16").Value = 0.08529816951135136
$ws.Range("P16").Value = 0.05852960232096958
$ws.Range("Q16").Value = 0.82882339264
$ws.Range("R16").Value = 4.97294035584
$ws.Range("S16").Value = 0.0005094751559339606
$ws.Range("T16").Value = 0.0003763841349621497
